$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "STEPS:" list below the existing content -----------------
# Shared-string insertion order matters (it drives the <sst> index each new
# string lands at), so the cells are populated in the same order the author's
# saved file shows them in xl/sharedStrings.xml: STEPS, 2), 3), 1), 4), 5).
$ws.Range("A62").Value = "STEPS:"
$ws.Range("A64").Value = "2) Spin Spinner"
$ws.Range("A65").Value = "3) Move specified number of spaces"
$ws.Range("A63").Value = "1) Initialise Players,"
$ws.Range("A66").Value = "4) Complete action indicated by space"
$ws.Range("A67").Value = "5) Next player"

# The numbered steps (rows 63-67) pick up a shared style (column A+B, except
# the last row which only has column A). Toggling ShrinkToFit at its default
# value stamps the cells with a fresh "applyAlignment" cell format without
# altering how they actually look, matching the new style added to the
# workbook.
$ws.Range("A63:B66").ShrinkToFit = $false
$ws.Range("A67").ShrinkToFit = $false

# --- Restore the view to where the author left it --------------------------
# Scroll the window so row 51 is at the top, then put the selection on B70
# (matches the saved sheetView's topLeftCell/selection).
$excel.ActiveWindow.ScrollRow = 51
$ws.Range("B70").Select()
